$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 57: J,K,L 7 -> 2
$ws.Range("J57:L57").Value = 2

# Row 58: J,K,L 2 -> 9
$ws.Range("J58:L58").Value = 9

# Rows 59-66: G 10 -> 5
$ws.Range("G59:G66").Value = 5

# Row 65: J,K,L 7 -> 2
$ws.Range("J65:L65").Value = 2

# Row 66: J,K,L 2 -> 9
$ws.Range("J66:L66").Value = 9

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("K52").Select()
